# ADD results from server
# Update the A2 "co2" price values on each yearly sheet with the latest
# results received from the server. The 2035 sheet is left unchanged.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2040", "2045", "2050")
$newValues  = @(57, 195, 355, 355, 355)

for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $ws.Range("A2").Value = $newValues[$i]
}
